$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the 2022 column (M) mirroring the existing 2013-2021 series in D:L.
# Copy formats from the adjacent existing cells so the new cells reuse the
# same cell styles as their neighbours, then set the values.

# M4: header year 2022 - takes the same style as K4 (2020 header cell)
$ws.Range("K4").Copy()
$ws.Range("M4").PasteSpecial(-4122)
$ws.Range("M4").Value = 2022

# M5: "Small enterprises" 2022 value - takes the same style as L5 (2021 value cell)
$ws.Range("L5").Copy()
$ws.Range("M5").PasteSpecial(-4122)
$ws.Range("M5").Value = 2.2

# M6: "Medium-sized enterprises" 2022 value - takes the same style as L6 (2021 value cell)
$ws.Range("L6").Copy()
$ws.Range("M6").PasteSpecial(-4122)
$ws.Range("M6").Value = 1.2

# Update the saved selection to match the authored workbook.
$ws.Range("M10").Select()
